# BGUSO-93 WPFGen: Apply conditional visibility : In progress
# Adds a new "NotifyPropertyChanged" worksheet (WPF INotifyPropertyChanged
# codegen helper table) after "WPF_temp_test", mirroring the pattern used
# by the existing codegen sheets in this workbook.

$wb = $excel.ActiveWorkbook

# --- add the new sheet right after the last existing sheet (WPF_temp_test) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "NotifyPropertyChanged"

# --- hide helper columns B:J (kept around for the TRIM/FIND scratch work) ---
$ws.Range("B1:J1").EntireColumn.Hidden = $true

# --- property declaration text, first the 5 that make up the original
#     fill-down block (rows 3-7), then the header captions, then the
#     extra row (row 2) that was inserted above the block afterwards.
#     The order mirrors how the shared-string table was actually built. ---
$ws.Range("A3").Value = "        public string NameUkr { get; set; }"
$ws.Range("A4").Value = "        public CountryInfo JurisdictionCountry { get; set; }"
$ws.Range("A5").Value = "        public string CourtRegion { get; set; }"
$ws.Range("A6").Value = "        public string CourtID { get; set; }"
$ws.Range("A7").Value = "        public CourtInstanceType Instance { get; set; }"

$ws.Range("H1").Value = "Type"
$ws.Range("G1").Value = "Modifier"
$ws.Range("I1").Value = "PropName"
$ws.Range("J1").Value = "FieldName"
$ws.Range("K1").Value = "FieldDecl"
$ws.Range("L1").Value = "Accessor"

$ws.Range("A2").Value = "public string ShortTermRatingValueOther { get; set; }"

# --- header row styling ---
$ws.Range("G1:M1").Font.Bold = $true

# --- row 2 : standalone (non-shared) formulas ---
$ws.Range("B2").Formula = "=TRIM(A2)"
$ws.Range("C2").Value = 1
$ws.Range("D2").Formula = '=FIND(" ",$B2)'
$ws.Range("E2").Formula = '=FIND(" ",$B2,D2+1)'
$ws.Range("F2").Formula = '=FIND(" ",$B2,E2+1)'
$ws.Range("G2").Formula = "=TRIM(MID(`$B2,C2,D2-C2))"
$ws.Range("H2").Formula = "=TRIM(MID(`$B2,D2,E2-D2))"
$ws.Range("I2").Formula = "=TRIM(MID(`$B2,E2,F2-E2))"
$ws.Range("J2").Formula = '="_"&I2'
$ws.Range("K2").Formula = '="private " & H2 & " " & J2 & ";"'
$ws.Range("L2").Formula = '=G2& " " &H2& " " &I2 & " { get { return " & J2 & "; } set { " &J2 & " = value; OnPropertyChanged(" & CHAR(34) & I2 & CHAR(34) & "); } }"'

# --- rows 3-7 : filled together so the engine records them as one shared
#     formula group (matches the fill-handle behaviour that produced the
#     original workbook) ---
$ws.Range("C3:C7").Value = 1
$ws.Range("B3:B7").Formula = "=TRIM(A3)"
$ws.Range("D3:D7").Formula = '=FIND(" ",$B3)'
$ws.Range("E3:F7").Formula = '=FIND(" ",$B3,D3+1)'
$ws.Range("G3:G7").Formula = "=TRIM(MID(`$B3,C3,D3-C3))"
$ws.Range("H3:H7").Formula = "=TRIM(MID(`$B3,D3,E3-D3))"
$ws.Range("I3:I7").Formula = "=TRIM(MID(`$B3,E3,F3-E3))"
$ws.Range("J3:J7").Formula = '="_"&I3'
$ws.Range("K3:K7").Formula = '="private " & H3 & " " & J3 & ";"'
$ws.Range("L3:L7").Formula = '=G3& " " &H3& " " &I3 & " { get { return " & J3 & "; } set { " &J3 & " = value; OnPropertyChanged(" & CHAR(34) & I3 & CHAR(34) & "); } }"'

# --- page setup to match the other codegen sheets ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection / activation : the new sheet becomes the active tab ---
$ws.Range("L2").Select()
$ws.Activate()
